$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.005670666694641
$ws.Range("B1").Value = 2.115151166915894
$ws.Range("C1").Value = 6.611392974853516
$ws.Range("D1").Value = 1.805461764335632
$ws.Range("E1").Value = 1.369054794311523
